$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add column C header
$ws.Range("C1").Value = "yearGroup3"

# 2) Fill column C (3-year groupings), entered bottom-up (C19 -> C5) so that the
#    shared-string table append order matches the source order.
$ws.Range("C19").Value = "2015-2017"
$ws.Range("C18").Value = "2015-2017"
$ws.Range("C17").Value = "2015-2017"

$ws.Range("C16").Value = "2012-2014"
$ws.Range("C15").Value = "2012-2014"
$ws.Range("C14").Value = "2012-2014"

$ws.Range("C13").Value = "2009-2011"
$ws.Range("C12").Value = "2009-2011"
$ws.Range("C11").Value = "2009-2011"

$ws.Range("C10").Value = "2006-2008"
$ws.Range("C9").Value = "2006-2008"
$ws.Range("C8").Value = "2006-2008"

$ws.Range("C7").Value = "2003-2005"
$ws.Range("C6").Value = "2003-2005"
$ws.Range("C5").Value = "2003-2005"

$ws.Range("C4").Value = "2000-2002"
$ws.Range("C3").Value = "2000-2002"
$ws.Range("C2").Value = "2000-2002"

# 3) Rename column B header last (yGroup1 -> yearGroup5)
$ws.Range("B1").Value = "yearGroup5"

# 4) Right-align the whole used range (A1:C19)
$ws.Range("A1:C19").HorizontalAlignment = -4152

# 5) Style the new C1 header cell with a (non-bold) red font
$ws.Range("C1").Font.Name = "Calibri"
$ws.Range("C1").Font.Size = 10
$ws.Range("C1").Font.Color = 255
$ws.Range("C1").Font.Bold = $false

# 6) Column widths
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 18.28515625

# 7) Selection
$ws.Range("E7").Select()
